$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 43 (hunk 0)
$ws.Range("H43").Value = 3368.875
$ws.Range("I43").Value = 3064.4285
$ws.Range("J43").Value = 5500
$ws.Range("K43").Value = 3064.4285
$ws.Range("L43").Value = 5500
$ws.Range("M43").Value = -2995.4285
$ws.Range("N43").Value = -5638

# row 98 (hunk 1)
$ws.Range("H98").Value = 940
$ws.Range("I98").Value = 940
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 940
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 558
$ws.Range("N98").ClearContents()

# row 100 (hunk 2)
$ws.Range("H100").Value = 2206.5881
$ws.Range("I100").Value = 2065.5715
$ws.Range("J100").Value = 2864.6667
$ws.Range("K100").Value = 2065.5715
$ws.Range("L100").Value = 2864.6667
$ws.Range("M100").Value = -1524.5715
$ws.Range("N100").Value = -3946.6667

# row 122 (hunk 3)
$ws.Range("H122").Value = 940
$ws.Range("I122").Value = 940
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2820
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -370
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# row 54 (hunk 4)
$ws.Range("H54").Value = 35999
$ws.Range("I54").Value = 35999
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 35999
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -35230
$ws.Range("N54").ClearContents()

# row 63 (hunk 5)
$ws.Range("H63").Value = 7445.5
$ws.Range("I63").Value = 3594.1667
$ws.Range("J63").Value = 18999.5
$ws.Range("K63").Value = 3594.1667
$ws.Range("L63").Value = 18999.5
$ws.Range("M63").Value = -2908.1667
$ws.Range("N63").Value = -20371.5

# row 66 (hunk 6)
$ws.Range("H66").Value = 7445.5
$ws.Range("I66").Value = 3594.1667
$ws.Range("J66").Value = 18999.5
$ws.Range("K66").Value = 17970.8335
$ws.Range("L66").Value = 94997.5
$ws.Range("M66").Value = -14538.8335
$ws.Range("N66").Value = -101861.5

# row 88 (hunk 7)
$ws.Range("H88").Value = 1572.25
$ws.Range("I88").Value = 1040
$ws.Range("J88").Value = 1952.4286
$ws.Range("K88").Value = 1040
$ws.Range("L88").Value = 1952.4286
$ws.Range("M88").Value = -634
$ws.Range("N88").Value = -2764.4286

# row 91 (hunk 8)
$ws.Range("H91").Value = 1572.25
$ws.Range("I91").Value = 1040
$ws.Range("J91").Value = 1952.4286
$ws.Range("K91").Value = 1040
$ws.Range("L91").Value = 1952.4286
$ws.Range("M91").Value = 364
$ws.Range("N91").Value = -4760.4286

$ws = $wb.Worksheets.Item("BSM")
# row 20 (hunk 9)
$ws.Range("H20").Value = 2991.1667
$ws.Range("I20").Value = 1952.875
$ws.Range("J20").Value = 5067.75
$ws.Range("K20").Value = 1952.875
$ws.Range("L20").Value = 5067.75
$ws.Range("M20").Value = -1705.875
$ws.Range("N20").Value = -5561.75

# row 105 (hunk 10)
$ws.Range("H105").Value = 3050.0908
$ws.Range("I105").Value = 2353.652
$ws.Range("J105").Value = 4651.9
$ws.Range("K105").Value = 2353.652
$ws.Range("L105").Value = 4651.9
$ws.Range("M105").Value = -606.652
$ws.Range("N105").Value = -8145.9

$ws = $wb.Worksheets.Item("CRP")
# row 16 (hunk 11)
$ws.Range("H16").Value = 2615
$ws.Range("I16").Value = 2615
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2615
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2328
$ws.Range("N16").ClearContents()

# row 31 (hunk 12)
$ws.Range("H31").Value = 1567.4445
$ws.Range("I31").Value = 1575.875
$ws.Range("J31").Value = 1500
$ws.Range("K31").Value = 1575.875
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = -1280.875
$ws.Range("N31").Value = -2090

# row 34 (hunk 13)
$ws.Range("H34").Value = 1567.4445
$ws.Range("I34").Value = 1575.875
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 1575.875
$ws.Range("L34").Value = 1500
$ws.Range("M34").Value = -1373.875
$ws.Range("N34").Value = -1904

# row 58 (hunk 14)
$ws.Range("H58").Value = 2494.9092
$ws.Range("I58").Value = 2816
$ws.Range("J58").Value = 1050
$ws.Range("K58").Value = 2816
$ws.Range("L58").Value = 1050
$ws.Range("M58").Value = -2613
$ws.Range("N58").Value = -1456

# row 105 (hunk 15)
$ws.Range("H105").Value = 1453.2
$ws.Range("I105").Value = 1355.5555
$ws.Range("J105").Value = 2332
$ws.Range("K105").Value = 1355.5555
$ws.Range("L105").Value = 2332
$ws.Range("M105").Value = 391.44450000000006
$ws.Range("N105").Value = -5826

# row 107 (hunk 16)
$ws.Range("H107").Value = 1497.1111
$ws.Range("I107").Value = 450
$ws.Range("J107").Value = 1628
$ws.Range("K107").Value = 450
$ws.Range("L107").Value = 1628
$ws.Range("M107").Value = 1470
$ws.Range("N107").Value = -5468

# row 113 (hunk 17)
$ws.Range("H113").Value = 2615
$ws.Range("I113").Value = 2615
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2615
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -445
$ws.Range("N113").ClearContents()

# row 122 (hunk 18)
$ws.Range("H122").Value = 3347.4443
$ws.Range("I122").Value = 3835.6667
$ws.Range("J122").Value = 2371
$ws.Range("K122").Value = 11507.000100000001
$ws.Range("L122").Value = 7113
$ws.Range("M122").Value = -9057.000100000001
$ws.Range("N122").Value = -12013

# row 136 (hunk 19)
$ws.Range("H136").Value = 2494.9092
$ws.Range("I136").Value = 2816
$ws.Range("J136").Value = 1050
$ws.Range("K136").Value = 8448
$ws.Range("L136").Value = 3150
$ws.Range("M136").Value = -5898
$ws.Range("N136").Value = -8250

# row 140 (hunk 20)
$ws.Range("H140").Value = 99999
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 99999
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 99999
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -110359

$ws = $wb.Worksheets.Item("CUL")
# row 4 (hunk 21)
$ws.Range("H4").Value = 1515692.8
$ws.Range("I4").Value = 16518.684
$ws.Range("J4").Value = 30000000
$ws.Range("K4").Value = 49556.052
$ws.Range("L4").Value = 90000000
$ws.Range("M4").Value = -49444.052
$ws.Range("N4").Value = -90000224

# row 121 (hunk 22)
$ws.Range("H121").Value = 1429404.1
$ws.Range("I121").Value = 900
$ws.Range("J121").Value = 1667488.1
$ws.Range("K121").Value = 2700
$ws.Range("L121").Value = 5002464.300000001
$ws.Range("M121").Value = -1390
$ws.Range("N121").Value = -5005084.300000001

# row 131 (hunk 23)
$ws.Range("H131").Value = 1974.7333
$ws.Range("I131").Value = 1565.6364
$ws.Range("J131").Value = 3099.75
$ws.Range("K131").Value = 4696.9092
$ws.Range("L131").Value = 9299.25
$ws.Range("M131").Value = 343.09079999999994
$ws.Range("N131").Value = -19379.25

# row 139 (hunk 24)
$ws.Range("H139").Value = 4256.7144
$ws.Range("I139").Value = 4132.8335
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 12398.500499999998
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = -7258.500499999998
$ws.Range("N139").Value = -25280

# row 140 (hunk 25)
$ws.Range("H140").Value = 835531.56
$ws.Range("I140").Value = 910852.6
$ws.Range("J140").Value = 7000
$ws.Range("K140").Value = 2732557.8
$ws.Range("L140").Value = 21000
$ws.Range("M140").Value = -2727377.8
$ws.Range("N140").Value = -31360

# row 141 (hunk 26)
$ws.Range("H141").Value = 5916.6665
$ws.Range("I141").Value = 6100
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 18300
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -13120
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
# row 70 (hunk 27)
$ws.Range("H70").Value = 6964.6665
$ws.Range("I70").Value = 7000
$ws.Range("J70").Value = 6894
$ws.Range("K70").Value = 7000
$ws.Range("L70").Value = 6894
$ws.Range("M70").Value = -6730
$ws.Range("N70").Value = -7434

# row 73 (hunk 28)
$ws.Range("H73").Value = 6964.6665
$ws.Range("I73").Value = 7000
$ws.Range("J73").Value = 6894
$ws.Range("K73").Value = 7000
$ws.Range("L73").Value = 6894
$ws.Range("M73").Value = -6064
$ws.Range("N73").Value = -8766

# row 122 (hunk 29)
$ws.Range("H122").Value = 1957.9474
$ws.Range("I122").Value = 1983.5
$ws.Range("J122").Value = 1498
$ws.Range("K122").Value = 5950.5
$ws.Range("L122").Value = 4494
$ws.Range("M122").Value = -3500.5
$ws.Range("N122").Value = -9394

$ws = $wb.Worksheets.Item("LTW")
# row 68 (hunk 30)
$ws.Range("H68").Value = 1400
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1400
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 1400
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -2898

# row 71 (hunk 31)
$ws.Range("H71").Value = 1400
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1400
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 7000
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -14488

# row 122 (hunk 32)
$ws.Range("H122").Value = 3387.111
$ws.Range("I122").Value = 3333.5
$ws.Range("J122").Value = 3574.75
$ws.Range("K122").Value = 10000.5
$ws.Range("L122").Value = 10724.25
$ws.Range("M122").Value = -7550.5
$ws.Range("N122").Value = -15624.25

$ws = $wb.Worksheets.Item("WVR")
# row 58 (hunk 33)
$ws.Range("H58").Value = 50499.5
$ws.Range("I58").Value = 6999
$ws.Range("J58").Value = 94000
$ws.Range("K58").Value = 6999
$ws.Range("L58").Value = 94000
$ws.Range("M58").Value = -6691
$ws.Range("N58").Value = -94616

# row 61 (hunk 34)
$ws.Range("H61").Value = 99999.836
$ws.Range("I61").Value = 171666.33
$ws.Range("J61").Value = 28333.334
$ws.Range("K61").Value = 171666.33
$ws.Range("L61").Value = 28333.334
$ws.Range("M61").Value = -171374.33
$ws.Range("N61").Value = -28917.334
